# Scheduled runner update: refresh market-price-derived profit figures
# (currentAveragePrice*, Leve price/profit columns) on several rows across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 231.84616
$ws.Range("I9").Value = 155.90909
$ws.Range("J9").Value = 649.5
$ws.Range("K9").Value = 155.90909
$ws.Range("L9").Value = 649.5
$ws.Range("M9").Value = 13.09091000000001
$ws.Range("H12").Value = 186.125
$ws.Range("I12").Value = 117.8
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 117.8
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = 52.2
$ws.Range("N12").Value = -640
$ws.Range("H76").Value = 3599.8
$ws.Range("I76").Value = 3599.8
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3599.8
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3284.8
$ws.Range("H79").Value = 3599.8
$ws.Range("I79").Value = 3599.8
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3599.8
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2507.8
$ws.Range("I86").Value = 2079.6667
$ws.Range("J86").Value = 2856.5557
$ws.Range("K86").Value = 2079.6667
$ws.Range("L86").Value = 2856.5557
$ws.Range("M86").Value = -956.6667000000002
$ws.Range("N86").Value = -5102.5557
$ws.Range("I89").Value = 2079.6667
$ws.Range("J89").Value = 2856.5557
$ws.Range("K89").Value = 10398.3335
$ws.Range("L89").Value = 14282.7785
$ws.Range("M89").Value = -4782.333500000001
$ws.Range("N89").Value = -25514.7785
$ws.Range("H116").Value = 203866.53
$ws.Range("I116").Value = 65073.5
$ws.Range("J116").Value = 330042
$ws.Range("K116").Value = 65073.5
$ws.Range("L116").Value = 330042
$ws.Range("M116").Value = -61631.5
$ws.Range("N116").Value = -336926
$ws.Range("H141").Value = 1935.8182
$ws.Range("I141").Value = 1749.25
$ws.Range("J141").Value = 2433.3333
$ws.Range("K141").Value = 5247.75
$ws.Range("L141").Value = 7299.999899999999
$ws.Range("M141").Value = -67.75
$ws.Range("N141").Value = -17659.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 300.125
$ws.Range("I5").Value = 175.25
$ws.Range("J5").Value = 425
$ws.Range("K5").Value = 175.25
$ws.Range("L5").Value = 425
$ws.Range("M5").Value = -63.25
$ws.Range("N5").Value = -649
$ws.Range("H61").Value = 11111761
$ws.Range("I61").Value = 11111761
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 11111761
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -11111549
$ws.Range("H74").Value = 3435.2917
$ws.Range("I74").Value = 2128.7896
$ws.Range("J74").Value = 8400
$ws.Range("K74").Value = 2128.7896
$ws.Range("L74").Value = 8400
$ws.Range("M74").Value = -1254.7896
$ws.Range("H77").Value = 3435.2917
$ws.Range("I77").Value = 2128.7896
$ws.Range("J77").Value = 8400
$ws.Range("K77").Value = 10643.948
$ws.Range("L77").Value = 42000
$ws.Range("M77").Value = -6275.948
$ws.Range("H132").Value = 1826850.9
$ws.Range("I132").Value = 2158278.2
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6474834.600000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6472304.600000001
$ws.Range("H136").Value = 11111761
$ws.Range("I136").Value = 11111761
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 33335283
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -33332733

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 300.125
$ws.Range("I4").Value = 175.25
$ws.Range("J4").Value = 425
$ws.Range("K4").Value = 175.25
$ws.Range("L4").Value = 425
$ws.Range("M4").Value = -60.25
$ws.Range("N4").Value = -655
$ws.Range("H134").Value = 2583995
$ws.Range("I134").Value = 3405788.8
$ws.Range("J134").Value = 1145855.8
$ws.Range("K134").Value = 10217366.4
$ws.Range("L134").Value = 3437567.4
$ws.Range("M134").Value = -10214831.4
$ws.Range("N134").Value = -3442637.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 240.22728
$ws.Range("I7").Value = 138.64285
$ws.Range("J7").Value = 418
$ws.Range("K7").Value = 138.64285
$ws.Range("L7").Value = 418
$ws.Range("M7").Value = -25.64285000000001
$ws.Range("N7").Value = -644
$ws.Range("H15").Value = 9021.111000000001
$ws.Range("I15").Value = 200
$ws.Range("J15").Value = 26663.334
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = 26663.334
$ws.Range("M15").Value = -30
$ws.Range("N15").Value = -27003.334
$ws.Range("H123").Value = 85000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 85000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 85000
$ws.Range("N123").Value = -94800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 285.2
$ws.Range("I8").Value = 285.2
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 855.5999999999999
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -716.5999999999999
$ws.Range("H75").Value = 5070.375
$ws.Range("I75").Value = 1815.5
$ws.Range("J75").Value = 5535.357
$ws.Range("K75").Value = 5446.5
$ws.Range("L75").Value = 16606.071
$ws.Range("M75").Value = -4448.5
$ws.Range("H78").Value = 5070.375
$ws.Range("I78").Value = 1815.5
$ws.Range("J78").Value = 5535.357
$ws.Range("K78").Value = 16339.5
$ws.Range("L78").Value = 49818.213
$ws.Range("M78").Value = -11347.5
$ws.Range("H131").Value = 16304.6
$ws.Range("I131").Value = 1130
$ws.Range("J131").Value = 23891.9
$ws.Range("K131").Value = 3390
$ws.Range("L131").Value = 71675.70000000001
$ws.Range("M131").Value = 1650
$ws.Range("N131").Value = -81755.70000000001
$ws.Range("H132").Value = 988.5
$ws.Range("I132").Value = 988.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8896.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6366.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 4264.4
$ws.Range("I29").Value = 1437.6666
$ws.Range("J29").Value = 8504.5
$ws.Range("K29").Value = 1437.6666
$ws.Range("L29").Value = 8504.5
$ws.Range("M29").Value = -1147.6666
$ws.Range("N29").Value = -9084.5
$ws.Range("H107").Value = 44079.777
$ws.Range("I107").Value = 92086.75
$ws.Range("J107").Value = 5674.2
$ws.Range("K107").Value = 92086.75
$ws.Range("L107").Value = 5674.2
$ws.Range("M107").Value = -90166.75
$ws.Range("N107").Value = -9514.200000000001
$ws.Range("H122").Value = 7845.4443
$ws.Range("I122").Value = 4248.8184
$ws.Range("J122").Value = 13497.286
$ws.Range("K122").Value = 12746.4552
$ws.Range("L122").Value = 40491.858
$ws.Range("M122").Value = -10296.4552
$ws.Range("N122").Value = -45391.858
$ws.Range("H126").Value = 1669586.2
$ws.Range("I126").Value = 3335633.2
$ws.Range("J126").Value = 3539.4
$ws.Range("K126").Value = 10006899.6
$ws.Range("L126").Value = 10618.2
$ws.Range("M126").Value = -10004429.6
$ws.Range("N126").Value = -15558.2
$ws.Range("H132").Value = 1214696
$ws.Range("I132").Value = 1730309.1
$ws.Range("J132").Value = 11598.667
$ws.Range("K132").Value = 5190927.300000001
$ws.Range("L132").Value = 34796.001
$ws.Range("M132").Value = -5188397.300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 25000
$ws.Range("I29").Value = 15000
$ws.Range("J29").Value = 30000
$ws.Range("K29").Value = 15000
$ws.Range("L29").Value = 30000
$ws.Range("M29").Value = -14705
$ws.Range("N29").Value = -30590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 77125630
$ws.Range("I4").Value = 853332
$ws.Range("J4").Value = 100007320
$ws.Range("K4").Value = 853332
$ws.Range("L4").Value = 100007320
$ws.Range("M4").Value = -853219
$ws.Range("N4").Value = -100007546
$ws.Range("H132").Value = 5034765
$ws.Range("I132").Value = 5298831.5
$ws.Range("J132").Value = 17500
$ws.Range("K132").Value = 15896494.5
$ws.Range("L132").Value = 52500
$ws.Range("M132").Value = -15893964.5
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

Write-Host "Updated profit figures on ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
